$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking "Price" values in column D must stay text (the sheet stores them
# as inline strings). Pre-formatting the cell as Text ("@") before assignment keeps
# COM from coercing the string into a real number.
$priceCells = @(
    @{ Ref = "D2"; Value = "242.87" }
    @{ Ref = "D3"; Value = "23.25" }
    @{ Ref = "D4"; Value = "5.557" }
    @{ Ref = "D5"; Value = "0.05872" }
    @{ Ref = "D6"; Value = "3.406" }
    @{ Ref = "D8"; Value = "1.319" }
    @{ Ref = "D9"; Value = "0.7980" }
    @{ Ref = "D10"; Value = "0.1462" }
    @{ Ref = "D11"; Value = "0.07606" }
    @{ Ref = "D12"; Value = "0.03255" }
    @{ Ref = "D13"; Value = "0.02978" }
    @{ Ref = "D14"; Value = "0.09238" }
    @{ Ref = "D15"; Value = "0.001669" }
    @{ Ref = "D16"; Value = "3.412" }
    @{ Ref = "D17"; Value = "0.04749" }
    @{ Ref = "D18"; Value = "0.0005995" }
    @{ Ref = "D19"; Value = "0.006244" }
    @{ Ref = "D20"; Value = "0.001068" }
    @{ Ref = "D23"; Value = "3.696" }
    @{ Ref = "D24"; Value = "2.209" }
    @{ Ref = "D25"; Value = "0.3335" }
    @{ Ref = "D26"; Value = "0.1252" }
    @{ Ref = "D27"; Value = "0.0004003" }
    @{ Ref = "D40"; Value = "0.04319" }
    @{ Ref = "D41"; Value = "0.007142" }
    @{ Ref = "D42"; Value = "0.003194" }
    @{ Ref = "D43"; Value = "0.1053" }
    @{ Ref = "D44"; Value = "0.008795" }
    @{ Ref = "D46"; Value = "0.00005724" }
    @{ Ref = "D48"; Value = "0.7858" }
    @{ Ref = "D49"; Value = "0.1054" }
    @{ Ref = "D50"; Value = "0.00002102" }
)
foreach ($item in $priceCells) {
    $rng = $ws.Range($item.Ref)
    $rng.NumberFormat = "@"
    $rng.Value = $item.Value
}

# Plain text fields (coin name, link, volume label) - assign directly.
$ws.Range("E18").Value = "17OneONE"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("E49").Value = "48BOLOBOLOBestin24h"
